$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Reorder names in the "Recorded By" column (G) for specific session rows.
# Each entry maps a row number to its new comma-separated value.
$updates = @{
    2 = "backup@backdoor.com, System, system"
    4 = "backup@backdoor.com, System"
    5 = "backup@backdoor.com, System"
    7 = "System, admin@admin.com"
    8 = "backup@backdoor.com, System"
    11 = "System, dnasr281@gmail.com"
    17 = "System, dnasr281@gmail.com"
    28 = "backup@backdoor.com, System, system"
    30 = "backup@backdoor.com, System"
    31 = "backup@backdoor.com, System"
    33 = "System, admin@admin.com"
    34 = "backup@backdoor.com, System"
    37 = "System, dnasr281@gmail.com"
    43 = "System, dnasr281@gmail.com"
    54 = "backup@backdoor.com, System, system"
    56 = "backup@backdoor.com, System"
    57 = "backup@backdoor.com, System"
    59 = "System, admin@admin.com"
    60 = "backup@backdoor.com, System"
    63 = "System, dnasr281@gmail.com"
    69 = "System, dnasr281@gmail.com"
    80 = "backup@backdoor.com, System"
    81 = "backup@backdoor.com, System"
    82 = "backup@backdoor.com, System"
    93 = "System, dnasr281@gmail.com"
    94 = "System, dnasr281@gmail.com"
    96 = "System, dnasr281@gmail.com"
    106 = "backup@backdoor.com, System"
    107 = "backup@backdoor.com, System"
    108 = "backup@backdoor.com, System"
    119 = "System, dnasr281@gmail.com"
    120 = "System, dnasr281@gmail.com"
    122 = "System, dnasr281@gmail.com"
    132 = "backup@backdoor.com, System"
    133 = "backup@backdoor.com, System"
    134 = "backup@backdoor.com, System"
    145 = "System, dnasr281@gmail.com"
    146 = "System, dnasr281@gmail.com"
    148 = "System, dnasr281@gmail.com"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item([int]$row, 7).Value = $updates[$row]
}
